$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at O:P, shifting the old O..U block
# (Extracted Objects .. Result String) two columns to the right
# (Q..W), matching the new "Correct Pred Predicates Parents/Related"
# columns being added before "Extracted Objects".
$ws.Range("O1:P1").EntireColumn.Insert()

# Rename the existing predicate-parent/related columns.
$ws.Range("M1").Value = "Detected Predicates Doc Parent"
$ws.Range("N1").Value = "Detected Predicates Doc Related"

# Header text for the two newly inserted columns.
$ws.Range("O1").Value = "Correct Pred Predicates Parents"
$ws.Range("P1").Value = "Correct Pred Predicates Related"

# New data values for the two newly inserted columns, per row.
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 2

$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1

$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1

$ws.Range("O5").Value = 2
$ws.Range("P5").Value = 2

$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 2
